$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new cell values in the same order the strings were first
# --- authored, so the shared-string table builds up in the same order.
$ws.Range("B30").Value = "max length "
$ws.Range("A31").Value = " Current UPB"
$ws.Range("A32").Value = "Loan Age maxlength"
$ws.Range("A33").Value = "Months to Legal Maturity"
$ws.Range("A34").Value = "Current Interest Rate"
$ws.Range("A35").Value = "Current Deferred UPB"
$ws.Range("A29").Value = "Testing for  Max Length"
$ws.Range("C30").Value = "Result"
$ws.Range("C31").Value = "UI doesn't allow User to enter"
$ws.Range("C32").Value = "UI doesn't allow User to enter"
$ws.Range("C33").Value = "UI doesn't allow User to enter"
$ws.Range("C34").Value = "UI doesn't allow User to enter"
$ws.Range("C35").Value = "UI doesn't allow User to enter"

$ws.Range("B31").Value = 9
$ws.Range("B32").Value = 3
$ws.Range("B33").Value = 3
$ws.Range("B34").Value = 5
$ws.Range("B35").Value = 9

# --- Formatting ---
# A29 ("Testing for  Max Length") reuses the same section-title look as
# A1/A13/A21 (fillId 33 style) - copy that formatting across.
$ws.Range("A13").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# B30 / C30 are the little table's header cells - bold text, no fill.
$ws.Range("B30:C30").Font.Bold = $true

# --- View changes ---
$excel.ActiveWindow.Zoom = 69
$ws.Range("A3").Select()
